# Fixed bug in altCount causing interval to drop to 0, leading to timer
# endless loop. Adds a new "Sheet2" with a cycleCount/altCount/interval/
# timeDir trace table, and makes Sheet2 the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row (cycleCount/altCount/interval first; timeDir's header and the
# round-0 comment are entered further below, matching shared-string order).
$ws2.Range("A1").Value = "cycleCount"
$ws2.Range("B1").Value = "altCount"
$ws2.Range("C1").Value = "interval"

# Row 2 is the seed row (no formula for interval here).
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 200
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = "Don't modify interval at beginning of round 0"
$ws2.Range("D1").Value = "timeDir"

# cycleCount / altCount / timeDir values for rows 3..26 (A/B/D columns).
$cycleCount = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1)
$altCount   = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,0,1,2,3,4,5)
$timeDir    = @(1,1,1,1,1,1,1,1,-1,-1,-1,-1,-1,-1,-1,-1,-1,1,1,1,1,1,1,1)

for ($i = 0; $i -lt 24; $i++) {
    $row = 3 + $i
    $ws2.Cells.Item($row, 1).Value = $cycleCount[$i]
    $ws2.Cells.Item($row, 2).Value = $altCount[$i]
    $ws2.Cells.Item($row, 4).Value = $timeDir[$i]
}

# interval column: C3 adds onto C2/D2, every other row continues the chain.
$ws2.Range("C3").Formula = "=C2+(D2*100)"
for ($row = 4; $row -le 26; $row++) {
    $prev = $row - 1
    $ws2.Range("C$row").Formula = "=C$prev+(D$prev*100)"
}

# View state: Sheet1's own selection stays put (D3), but Sheet2 becomes
# the active/selected tab with A22 selected, so select on Sheet1 first.
$ws1.Range("D3").Select()
$ws2.Activate()
$ws2.Range("A22").Select()
